$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision has been recorded with no movie selected for Friday.`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selected for Friday.`n"
$ws.Range("D3").Value = "no_decision, "
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision about which movie to play on Friday has resulted in no consensus.`n"
$ws.Range("D4").Value = "no_decision, "
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie was acquired for this Friday.`n"
$ws.Range("D5").Value = "no_decision, "
$ws.Range("C6").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for the movie `"Barbie.`"`n"
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision has been recorded as no choice of movie for Friday.`n"
$ws.Range("D8").Value = "no_decision, "
$ws.Range("C9").Value = "MSG: None`n`nMSG: The committee did not reach a consensus regarding which movie to show on Friday, resulting in no decision being made.`n"
$ws.Range("D9").Value = "no_decision, "
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision about what movie to show on Friday ended without any agreement.`n"
$ws.Range("D10").Value = "no_decision, "
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision has been recorded as no consensus was reached regarding the movie to be shown on Friday.`n"
$ws.Range("D12").Value = "no_decision, "
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision.`"`n"
$ws.Range("D13").Value = "no_decision, "
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been recorded successfully.`n"
$ws.Range("D14").Value = "both_movies, "
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no agreement was reached regarding which movie to show on Friday.`n"
$ws.Range("D15").Value = "no_decision, "
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision outcome is recorded as no decision being made regarding the movie to be shown on Friday.`n"
$ws.Range("D16").Value = "no_decision, "
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision about the movie for Friday could not be made, leading to no acquisition.`n"
$ws.Range("D17").Value = "no_decision, "
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie was selected during the conversation.`n"
$ws.Range("D18").Value = "no_decision, "
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision has been recorded with no movie selected for Friday.`n"
$ws.Range("D19").Value = "no_decision, "
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision has been recorded with no choice on the movie for Friday.`n"
$ws.Range("D20").Value = "no_decision, "
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Barbie.`"`n"
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie selection.`n"
$ws.Range("D22").Value = "no_decision, "
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been confirmed.`n"
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected.`n"
$ws.Range("D24").Value = "no_decision, "
$ws.Range("C25").Value = "MSG: None`n`nMSG: The rights for both movies have been successfully acquired.`n"
$ws.Range("D25").Value = "both_movies, "
$ws.Range("C26").Value = "MSG: None`n`nMSG: The committee did not come to a decision about what movie to show on Friday.`n"
$ws.Range("D26").Value = "no_decision, "
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision has been recorded, reflecting that no agreement was reached on which movie to show on Friday.`n"
$ws.Range("D27").Value = "no_decision, "
$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision has been recorded as no selection for Friday's movie.`n"
$ws.Range("D28").Value = "no_decision, "
$ws.Range("C29").Value = "MSG: None`n`nMSG: The decision was made not to select a movie for Friday, as there was no consensus between the committee members.`n"
$ws.Range("D29").Value = "no_decision, "
$ws.Range("C30").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been finalized.`n"
$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has not been made, so I have called the no_decision function.`n"
$ws.Range("D32").Value = "no_decision, "
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie`" as the movie to be shown on Friday.`n"
$ws.Range("C34").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C35").Value = "MSG: None`n`nMSG: It seems no decision about Friday’s movie was reached.`n"
$ws.Range("D35").Value = "no_decision, "
$ws.Range("C36").Value = "MSG: None`n`nMSG: I have recorded the decision to acquire the rights for `"Oppenheimer`" as the selected movie to show on Friday.`n"
$ws.Range("C37").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision.`"`n"
$ws.Range("D37").Value = "no_decision, "
$ws.Range("C38").Value = "MSG: None`n`nMSG: No movie was selected in this meeting.`n"
$ws.Range("D38").Value = "no_decision, "
$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Barbie`" for Friday's movie.`n"
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision.`"`n"
$ws.Range("D40").Value = "no_decision, "
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision about Friday's movie was not finalized, so no agreement was reached.`n"
$ws.Range("D41").Value = "no_decision, "
